$wb = $excel.ActiveWorkbook

# The "adductName" column ([M-H]- values) is removed from the Corrected,
# Normalized, and PoolAfterDF sheets. On "Corrected"/"Normalized" that is
# column C; on "PoolAfterDF" it is column B.

$wsCorrected = $wb.Worksheets.Item("Corrected")
$wsCorrected.Columns.Item(3).Delete()

$wsNormalized = $wb.Worksheets.Item("Normalized")
$wsNormalized.Columns.Item(3).Delete()

$wsPoolAfterDF = $wb.Worksheets.Item("PoolAfterDF")
$wsPoolAfterDF.Columns.Item(2).Delete()
